$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename score headers (row 1)
$ws.Range("D1").Value = "HW_Score"
$ws.Range("E1").Value = "Q_Score"
$ws.Range("F1").Value = "Midterm_Score"
$ws.Range("G1").Value = "Final_Score"
$ws.Range("H1").Value = "HMN_Score"
$ws.Range("I1").Value = "QN_Score"
$ws.Range("J1").Value = "MidtermN_Score"
$ws.Range("K1").Value = "FinalN_Score"

# Set column widths (B and C)
$ws.Columns.Item(2).ColumnWidth = 15.625
$ws.Columns.Item(3).ColumnWidth = 39.0625

# Update score values (U5 student additions / corrections)
$ws.Range("D2").Value = 20.0
$ws.Range("G2").Value = 0.0
$ws.Range("H2").Value = 20.0
$ws.Range("K2").Value = 0.0

$ws.Range("G4").Value = 0.0
$ws.Range("K4").Value = 0.0

$ws.Range("G5").Value = 0.0
$ws.Range("K5").Value = 0.0

$ws.Range("G6").Value = 0.0
$ws.Range("K6").Value = 0.0

$ws.Range("D10").Value = 0.0
$ws.Range("G10").Value = 0.0
$ws.Range("H10").Value = 0.0
$ws.Range("K10").Value = 0.0

$ws.Range("G14").Value = 0.0
$ws.Range("K14").Value = 0.0

$ws.Range("G38").Value = 0.0
$ws.Range("K38").Value = 0.0
